# Susan David Bio (zh_CN) - reflow into single runs per paragraph,
# add first-line indent, stamp a _GoBack bookmark, and drop the
# trailing blank paragraph at the end of the body.

$d = $word.ActiveDocument

# --- 1. Collapse each multi-run bio paragraph down to one run ------------
# A Range.Text assignment that is identical to the existing text is a
# no-op, so stage through a placeholder string first - that forces Word
# to drop the old runs, then we write the real (unchanged) text back into
# the single remaining run.

# Paragraph 1: the opening bio paragraph (currently 7 runs).
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
[void]$r1.MoveEnd(1, -1)            # exclude the paragraph mark
$p1Text = $r1.Text
$r1.Text = "PLACEHOLDER1"
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
[void]$r1.MoveEnd(1, -1)
$r1.Text = $p1Text

# Paragraph 3: the awards paragraph (currently 3 runs).
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
[void]$r3.MoveEnd(1, -1)
$p3Text = $r3.Text
$r3.Text = "PLACEHOLDER3"
$p3 = $d.Paragraphs.Item(3)
$r3 = $p3.Range
[void]$r3.MoveEnd(1, -1)
$r3.Text = $p3Text

# --- 2. First-line indent on the three text-bearing paragraphs -----------
$d.Paragraphs.Item(1).Format.FirstLineIndent = 36   # 36pt = 720 twips
$d.Paragraphs.Item(3).Format.FirstLineIndent = 36
$d.Paragraphs.Item(5).Format.FirstLineIndent = 36

# --- 3. Stamp a _GoBack bookmark at the start of the final bio paragraph -
$p5 = $d.Paragraphs.Item(5)
$goBackRange = $p5.Range.Duplicate
$goBackRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# --- 4. Drop the trailing empty paragraph at the very end of the body ----
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$prevPara = $d.Paragraphs.Item($d.Paragraphs.Count - 1)
$trailRange = $d.Range($prevPara.Range.End - 1, $lastPara.Range.End)
$trailRange.Delete()
